$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like numbers/dates to Excel's auto-detection;
# force them to Text format first so they are stored as literal strings,
# then restore the default "Normal" style so no stray formatting is left behind.
$textCells = @("D5", "D10", "D13", "D18", "D20", "D24", "D28", "D29", "D32", "D34", "D35", "D40", "D43", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.746.22'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.602.98'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '211.92'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = '19.67'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '1.828.20'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.07'
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.585.28'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').Value = '209.60'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '7.14'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('E21').Value = '  +0.28%  '
$ws.Range('E22').Value = '  -4.74%  '
$ws.Range('D24').Value = '143.74'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = '15.37'
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('D29').Value = '0.0507'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D32').Value = '2.96'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = '1.288.56'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '1.24'
$ws.Range('E34').Value = '  +19.16%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  -4.57%  '
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = '5.44'
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').Value = '62.83'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').Value = '1.739.80'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').Value = '90.50'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.102'
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.0513'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Value = '5.96'
$ws.Range('E49').Value = '  +16.38%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.57'
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.25%  '

foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
